# Add a new row (row 13) for year 2021 to the worksheet, following the
# same pattern as the existing rows (e.g. row 12 / 2020).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 13 by copying the formatting/shape of row 12 (the previous
# year's row) down one row, then overwrite with the 2021 data. This keeps
# the year-label style (bold/centered, bordered) and the "present but
# blank" cells in columns B and E consistent with the rest of the sheet.
$ws.Range("A12:G12").Copy($ws.Range("A13:G13"))

$ws.Range("A13").Value = "2021年"
$ws.Range("C13").Value = 123687.09
$ws.Range("D13").Value = 9518
$ws.Range("F13").Value = 29625.95
$ws.Range("G13").Value = 18956.36
